$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 192
    3  = 433
    4  = 12618
    5  = 1283
    6  = 152
    7  = 34
    14 = 127
    15 = 44
    17 = 4407
    23 = 80
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
